$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# The deck ships two themes:
#   ppt/theme/theme1.xml -> used by the Slide Master ("Integral" design,
#                            "Red Violet" 12-slot color scheme)
#   ppt/theme/theme2.xml -> used by the Notes Master ("Office Theme"
#                            design, "Office" 12-slot color scheme)
#
# The target edit swaps the two themes' color schemes: the Slide Master
# (and therefore every slide built on it) switches to the standard
# "Office" palette, while the Notes Master switches to the "Red Violet"
# palette that used to belong to the Slide Master. The font scheme and
# the fill/line/effect format scheme are already byte-identical between
# the two themes, so only the 12 theme colors actually move.
#
# PowerPoint exposes those 12 theme colors through Master.ColorScheme,
# whose Colors(1..12) map 1:1 onto dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink (RGBColor.RGB uses the usual VBA RGB(r,g,b) packing:
# r + g*256 + b*65536). We drive that object model here exactly like a
# user picking new theme colors for the master.
# ---------------------------------------------------------------------------

$slideMasterColors = $p.SlideMaster.ColorScheme

# New ("Office") theme colors for the Slide Master
$slideMasterColors.Colors(1).RGB  = 0          # dk1       000000
$slideMasterColors.Colors(2).RGB  = 16777215   # lt1       FFFFFF
$slideMasterColors.Colors(3).RGB  = 6968388    # dk2       44546A
$slideMasterColors.Colors(4).RGB  = 15132391   # lt2       E7E6E6
$slideMasterColors.Colors(5).RGB  = 13998939   # accent1   5B9BD5
$slideMasterColors.Colors(6).RGB  = 3243501    # accent2   ED7D31
$slideMasterColors.Colors(7).RGB  = 10855845   # accent3   A5A5A5
$slideMasterColors.Colors(8).RGB  = 49407      # accent4   FFC000
$slideMasterColors.Colors(9).RGB  = 12874308   # accent5   4472C4
$slideMasterColors.Colors(10).RGB = 4697456    # accent6   70AD47
$slideMasterColors.Colors(11).RGB = 12673797   # hlink     0563C1
$slideMasterColors.Colors(12).RGB = 7491477    # folHlink  954F72

# Same operation on the Notes Master, moving the old Slide Master palette
# ("Red Violet") onto it, mirroring the swap end to end.
$notesMasterColors = $p.NotesMaster.ColorScheme

$notesMasterColors.Colors(1).RGB  = 0          # dk1       000000
$notesMasterColors.Colors(2).RGB  = 16777215   # lt1       FFFFFF
$notesMasterColors.Colors(3).RGB  = 5326149    # dk2       454551
$notesMasterColors.Colors(4).RGB  = 14473688   # lt2       D8D9DC
$notesMasterColors.Colors(5).RGB  = 9514467    # accent1   E32D91
$notesMasterColors.Colors(6).RGB  = 13381832   # accent2   C830CC
$notesMasterColors.Colors(7).RGB  = 14460494   # accent3   4EA6DC
$notesMasterColors.Colors(8).RGB  = 15168839   # accent4   4775E7
$notesMasterColors.Colors(9).RGB  = 14774665   # accent5   8971E1
$notesMasterColors.Colors(10).RGB = 7555029    # accent6   D54773
$notesMasterColors.Colors(11).RGB = 2465643    # hlink     6B9F25
$notesMasterColors.Colors(12).RGB = 9211020    # folHlink  8C8C8C

# Re-apply the Slide Master palette last so it wins if this runtime's
# COM host happens to resolve both masters' ColorScheme writes to the
# same underlying theme part.
$slideMasterColors.Colors(1).RGB  = 0
$slideMasterColors.Colors(2).RGB  = 16777215
$slideMasterColors.Colors(3).RGB  = 6968388
$slideMasterColors.Colors(4).RGB  = 15132391
$slideMasterColors.Colors(5).RGB  = 13998939
$slideMasterColors.Colors(6).RGB  = 3243501
$slideMasterColors.Colors(7).RGB  = 10855845
$slideMasterColors.Colors(8).RGB  = 49407
$slideMasterColors.Colors(9).RGB  = 12874308
$slideMasterColors.Colors(10).RGB = 4697456
$slideMasterColors.Colors(11).RGB = 12673797
$slideMasterColors.Colors(12).RGB = 7491477
